$d = $word.ActiveDocument

$replacements = @(
    @("69×78=", "91×81="),
    @("65×54=", "66×79="),
    @("57×58=", "67×56="),
    @("66×77=", "57×68="),
    @("50×72=", "59×49="),
    @("11×83=", "30×34="),
    @("64×79=", "54×28="),
    @("53×14=", "50×71="),
    @("42×65=", "14×59="),
    @("69×48=", "86×77="),
    @("57×67=", "59×25="),
    @("98×59=", "45×24="),
    @("43×43=", "59×45="),
    @("59×61=", "21×57="),
    @("99×67=", "65×88="),
    @("54×66=", "57×62="),
    @("39×96=", "60×50="),
    @("41×63=", "60×52="),
    @("52×62=", "13×69="),
    @("68×65=", "31×49="),
    @("68×43=", "63×68="),
    @("13×91=", "59×29="),
    @("82×82=", "19×50="),
    @("33×43=", "61×45="),
    @("62×93=", "82×86=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
